# LAU Report - Valuation
#
# Inserts a new "Residential Value Post-Capex" column into the LTC sheet
# (right before the existing "Vacant Value Pre-Capex" column becomes a
# three-line label), re-wires the AutoFilter / _FilterDatabase range to
# include the extra column, and makes sure the new column's subtotal
# formula (row 1) is filled in just like its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column at X. Everything that used to live in X..AF
# (values, formulas, formatting) shifts right to Y..AG automatically.
$ws.Columns("X").Insert()

# The old "Vacant Value Pre-Capex" header (column W) now wraps onto a
# third line.
$ws.Range("W3").Value = "Vacant`nValue`nPre-Capex"

# Give the brand-new column (X) its header.
$ws.Range("X3").Value = "Residential`nValue`nPost-Capex"

# Row 1 carries a SUBTOTAL() per column; fill it in for the new column,
# matching the pattern used by its neighbours.
$ws.Range("X1").Formula = "=SUBTOTAL(9,X4:X99999)"

# Re-apply the AutoFilter so it spans the newly widened table.
$ws.AutoFilterMode = $false
$ws.Range("A3:AG4").AutoFilter()

# Keep the _FilterDatabase defined name synced with the AutoFilter range.
$wb.Names.Item(1).RefersTo = "=LTC!`$A`$3:`$AG`$4"
